$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for Price column (D) values so numeric-looking
# strings like "208.56" are not auto-converted to numbers by Excel,
# and restore the default "Normal" style afterwards so no new
# cell style is left applied.
function Set-TextValue($cell, $text) {
    $cell.Style = "Normal"
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "26.121.15"
$ws.Range("E2").Value = "  -0.23%  "
Set-TextValue $ws.Range("D3") "1.664.55"
$ws.Range("E3").Value = "  -0.78%  "
Set-TextValue $ws.Range("D4") "1.003"
$ws.Range("E4").Value = "  -0.17%  "
Set-TextValue $ws.Range("D5") "208.56"
$ws.Range("E5").Value = "  -0.80%  "
Set-TextValue $ws.Range("D6") "0.5215"
$ws.Range("E6").Value = "  -1.41%  "
$ws.Range("E7").Value = "  -0.09%  "
Set-TextValue $ws.Range("D8") "0.2594"
$ws.Range("E8").Value = "  -3.25%  "
Set-TextValue $ws.Range("D9") "0.06323"
$ws.Range("E9").Value = "  +0.47%  "
Set-TextValue $ws.Range("D10") "21.00"
$ws.Range("E10").Value = "  -1.45%  "
Set-TextValue $ws.Range("D11") "0.07533"
Set-TextValue $ws.Range("D12") "1.664.14"
Set-TextValue $ws.Range("D13") "4.407"
$ws.Range("E13").Value = "  -1.74%  "
Set-TextValue $ws.Range("D14") "0.5377"
$ws.Range("E14").Value = "  -4.98%  "
Set-TextValue $ws.Range("D15") "0.0₅7972"
$ws.Range("E15").Value = "  -1.83%  "
Set-TextValue $ws.Range("D16") "66.06"
$ws.Range("E16").Value = "  -0.20%  "
Set-TextValue $ws.Range("D17") "26.145.17"
$ws.Range("E17").Value = "  -0.18%  "
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("E19").Value = "  -3.02%  "
$ws.Range("E20").Value = "  -0.73%  "
Set-TextValue $ws.Range("D21") "10.20"
$ws.Range("E21").Value = "  -3.07%  "
Set-TextValue $ws.Range("D22") "6.197"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("E23").Value = "  -0.16%  "
Set-TextValue $ws.Range("D24") "149.09"
$ws.Range("E24").Value = "  +0.66%  "
Set-TextValue $ws.Range("D25") "0.1224"
$ws.Range("E25").Value = "  -2.97%  "
Set-TextValue $ws.Range("D26") "7.409"
$ws.Range("E26").Value = "  -2.60%  "
Set-TextValue $ws.Range("D27") "15.67"
$ws.Range("E27").Value = "  -1.24%  "
Set-TextValue $ws.Range("D28") "0.06187"
$ws.Range("E28").Value = "  -4.21%  "
Set-TextValue $ws.Range("D29") "1.360"
$ws.Range("E29").Value = "  +1.35%  "
$ws.Range("E30").Value = "  -1.12%  "
Set-TextValue $ws.Range("D31") "3.471"
$ws.Range("E31").Value = "  -1.55%  "
Set-TextValue $ws.Range("D32") "3.396"
$ws.Range("E32").Value = "  -2.70%  "
Set-TextValue $ws.Range("D33") "1.631"
$ws.Range("E33").Value = "  -0.89%  "
Set-TextValue $ws.Range("D34") "0.9896"
$ws.Range("E34").Value = "  -1.58%  "
Set-TextValue $ws.Range("D35") "2.392"
$ws.Range("E35").Value = "  -0.96%  "
Set-TextValue $ws.Range("D36") "2.754"
$ws.Range("E36").Value = "  +1.53%  "
Set-TextValue $ws.Range("D37") "0.5889"
$ws.Range("E37").Value = "  -3.38%  "
Set-TextValue $ws.Range("D38") "1.106.70"
$ws.Range("E38").Value = "  +0.52%  "
Set-TextValue $ws.Range("D39") "0.01597"
$ws.Range("E39").Value = "  -0.85%  "
Set-TextValue $ws.Range("D40") "6.007"
$ws.Range("E40").Value = "  -2.50%  "
Set-TextValue $ws.Range("D41") "0.8475"
$ws.Range("E41").Value = "  -2.17%  "
Set-TextValue $ws.Range("D42") "1.004"
Set-TextValue $ws.Range("D43") "99.92"
$ws.Range("E43").Value = "  -0.16%  "
Set-TextValue $ws.Range("D44") "1.816.67"
$ws.Range("E44").Value = "  -0.64%  "
Set-TextValue $ws.Range("D45") "0.0₈108"
$ws.Range("E45").Value = "  -0.57%  "
Set-TextValue $ws.Range("D46") "55.19"
$ws.Range("E46").Value = "  -2.85%  "
Set-TextValue $ws.Range("D47") "1.004"
$ws.Range("E47").Value = "  -0.20%  "
Set-TextValue $ws.Range("D48") "8.052"
$ws.Range("E48").Value = "  +0.91%  "
Set-TextValue $ws.Range("D49") "0.05242"
$ws.Range("E49").Value = "  -0.54%  "
Set-TextValue $ws.Range("D50") "0.4253"
$ws.Range("E50").Value = "  -0.42%  "
Set-TextValue $ws.Range("D51") "5.870"
$ws.Range("E51").Value = "  -1.43%  "
